# Auto-generated cell updates applying the Sheets market-data refresh diff.
# For each affected sheet/cell, set the new value (or clear it when removed).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (50 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 140.33333
$ws.Range("I9").Value = 140.33333
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 140.33333
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 28.66667000000001
$ws.Range("N9").Value = $null
$ws.Range("H32").Value = 3441.818
$ws.Range("I32").Value = 3366.25
$ws.Range("J32").Value = 3485
$ws.Range("K32").Value = 3366.25
$ws.Range("L32").Value = 3485
$ws.Range("M32").Value = -3040.25
$ws.Range("N32").Value = -4137
$ws.Range("H33").Value = 6098.3335
$ws.Range("I33").Value = 7568.7144
$ws.Range("K33").Value = 7568.7144
$ws.Range("M33").Value = -7339.7144
$ws.Range("H43").Value = 3171.5454
$ws.Range("I43").Value = 2559.4
$ws.Range("J43").Value = 3681.6667
$ws.Range("K43").Value = 2559.4
$ws.Range("L43").Value = 3681.6667
$ws.Range("M43").Value = -2490.4
$ws.Range("N43").Value = -3819.6667
$ws.Range("H62").Value = 16842.766
$ws.Range("I62").Value = 30132.834
$ws.Range("J62").Value = 9593.637000000001
$ws.Range("K62").Value = 30132.834
$ws.Range("L62").Value = 9593.637000000001
$ws.Range("M62").Value = -29508.834
$ws.Range("N62").Value = -10841.637
$ws.Range("H65").Value = 16842.766
$ws.Range("I65").Value = 30132.834
$ws.Range("J65").Value = 9593.637000000001
$ws.Range("K65").Value = 150664.17
$ws.Range("L65").Value = 47968.185
$ws.Range("M65").Value = -147544.17
$ws.Range("N65").Value = -54208.185
$ws.Range("H94").Value = 6167.3335
$ws.Range("I94").Value = 6167.3335
$ws.Range("K94").Value = 6167.3335
$ws.Range("M94").Value = -5716.3335
$ws.Range("H100").Value = 1775.2593
$ws.Range("I100").Value = 1263.6666
$ws.Range("J100").Value = 3565.8333
$ws.Range("K100").Value = 1263.6666
$ws.Range("L100").Value = 3565.8333
$ws.Range("M100").Value = -722.6666
$ws.Range("N100").Value = -4647.8333

# --- Sheet: ARM (16 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 34999
$ws.Range("I7").Value = 34999
$ws.Range("K7").Value = 34999
$ws.Range("M7").Value = -34885
$ws.Range("H32").Value = 11395.926
$ws.Range("I32").Value = 7457.4736
$ws.Range("K32").Value = 7457.4736
$ws.Range("M32").Value = -7170.4736
$ws.Range("H74").Value = 203559.4
$ws.Range("I74").Value = 203559.4
$ws.Range("K74").Value = 203559.4
$ws.Range("M74").Value = -202685.4
$ws.Range("H77").Value = 203559.4
$ws.Range("I77").Value = 203559.4
$ws.Range("K77").Value = 1017797
$ws.Range("M77").Value = -1013429

# --- Sheet: BSM (12 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1251.7858
$ws.Range("I107").Value = 1263.4615
$ws.Range("K107").Value = 1263.4615
$ws.Range("M107").Value = 656.5385000000001
$ws.Range("H124").Value = 129999
$ws.Range("J124").Value = 129999
$ws.Range("L124").Value = 129999
$ws.Range("N124").Value = -139819
$ws.Range("H134").Value = 2022.807
$ws.Range("I134").Value = 1746.08
$ws.Range("K134").Value = 5238.24
$ws.Range("M134").Value = -2703.24

# --- Sheet: CRP (27 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2554.1035
$ws.Range("I31").Value = 2215.5186
$ws.Range("K31").Value = 2215.5186
$ws.Range("M31").Value = -1920.5186
$ws.Range("H34").Value = 2554.1035
$ws.Range("I34").Value = 2215.5186
$ws.Range("K34").Value = 2215.5186
$ws.Range("M34").Value = -2013.5186
$ws.Range("H55").Value = 24402.166
$ws.Range("I55").Value = 24402.166
$ws.Range("K55").Value = 24402.166
$ws.Range("M55").Value = -24087.166
$ws.Range("H86").Value = 45305.8
$ws.Range("I86").Value = 55158.125
$ws.Range("K86").Value = 55158.125
$ws.Range("M86").Value = -54035.125
$ws.Range("H89").Value = 45305.8
$ws.Range("I89").Value = 55158.125
$ws.Range("K89").Value = 275790.625
$ws.Range("M89").Value = -270174.625
$ws.Range("H134").Value = 24177.39
$ws.Range("I134").Value = 26749.023
$ws.Range("J134").Value = 3090
$ws.Range("K134").Value = 80247.069
$ws.Range("L134").Value = 9270
$ws.Range("M134").Value = -77712.069
$ws.Range("N134").Value = -14340

# --- Sheet: CUL (4 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1506277.9
$ws.Range("I4").Value = 1131937.5
$ws.Range("K4").Value = 3395812.5
$ws.Range("M4").Value = -3395700.5

# --- Sheet: GSM (19 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3692.2
$ws.Range("I80").Value = 3245.5
$ws.Range("K80").Value = 3245.5
$ws.Range("M80").Value = -2247.5
$ws.Range("H83").Value = 3692.2
$ws.Range("I83").Value = 3245.5
$ws.Range("K83").Value = 16227.5
$ws.Range("M83").Value = -11235.5
$ws.Range("H102").Value = 3014.4348
$ws.Range("I102").Value = 2502.0557
$ws.Range("K102").Value = 2502.0557
$ws.Range("M102").Value = -880.0556999999999
$ws.Range("H132").Value = 31020.543
$ws.Range("I132").Value = 32816.22
$ws.Range("J132").Value = 11866.667
$ws.Range("K132").Value = 98448.66
$ws.Range("L132").Value = 35600.001
$ws.Range("M132").Value = -95918.66
$ws.Range("N132").Value = -40660.001

# --- Sheet: LTW (31 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 71416.56
$ws.Range("J22").Value = 3088.4443
$ws.Range("L22").Value = 3088.4443
$ws.Range("N22").Value = -3678.4443
$ws.Range("H27").Value = 71416.56
$ws.Range("J27").Value = 3088.4443
$ws.Range("L27").Value = 3088.4443
$ws.Range("N27").Value = -3302.4443
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").Value = $null
$ws.Range("H68").Value = 3471
$ws.Range("J68").Value = 5000
$ws.Range("L68").Value = 5000
$ws.Range("N68").Value = -6498
$ws.Range("H71").Value = 3471
$ws.Range("J71").Value = 5000
$ws.Range("L71").Value = 25000
$ws.Range("N71").Value = -32488
$ws.Range("H93").Value = 1065.5264
$ws.Range("I93").Value = 802.875
$ws.Range("J93").Value = 2466.3333
$ws.Range("K93").Value = 802.875
$ws.Range("L93").Value = 2466.3333
$ws.Range("M93").Value = 445.125
$ws.Range("N93").Value = -4962.3333
$ws.Range("H140").Value = 69332.336
$ws.Range("J140").Value = 69332.336
$ws.Range("L140").Value = 69332.336
$ws.Range("N140").Value = -79692.336

# --- Sheet: WVR (59 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 25000
$ws.Range("J2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("N2").Value = -10224
$ws.Range("H3").Value = 7999.6
$ws.Range("I3").Value = 4999.5
$ws.Range("J3").Value = 9999.666999999999
$ws.Range("K3").Value = 4999.5
$ws.Range("L3").Value = 9999.666999999999
$ws.Range("M3").Value = -4885.5
$ws.Range("N3").Value = -10227.667
$ws.Range("H4").Value = 8616.333000000001
$ws.Range("I4").Value = 15732.667
$ws.Range("K4").Value = 15732.667
$ws.Range("M4").Value = -15619.667
$ws.Range("H6").Value = 19598.5
$ws.Range("J6").Value = 19598.5
$ws.Range("L6").Value = 19598.5
$ws.Range("N6").Value = -19828.5
$ws.Range("H10").Value = 205
$ws.Range("I10").Value = 205
$ws.Range("K10").Value = 205
$ws.Range("M10").Value = -36
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = $null
$ws.Range("H13").Value = 500
$ws.Range("I13").Value = 500
$ws.Range("K13").Value = 500
$ws.Range("M13").Value = -360
$ws.Range("H17").Value = 7950
$ws.Range("J17").Value = 7950
$ws.Range("L17").Value = 7950
$ws.Range("N17").Value = -8294
$ws.Range("H18").Value = 9999.5
$ws.Range("J18").Value = 9999.5
$ws.Range("L18").Value = 9999.5
$ws.Range("N18").Value = -10345.5
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = $null
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = $null
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = $null
$ws.Range("H132").Value = 36812.1
$ws.Range("I132").Value = 38119.68
$ws.Range("K132").Value = 114359.04
$ws.Range("M132").Value = -111829.04
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null
